$wb = $excel.ActiveWorkbook

# --- 1. Status text change: "Ready for handoff" -> "In Translation" ---
# This shared string is referenced by every Status/zh-cn/de-de status cell
# across all three sheets. Updating every cell that holds it collapses the
# old shared-string entry (it becomes unreferenced and is dropped on save)
# and replaces it everywhere, matching the source diff which edits the
# shared string itself.

$wsOverview = $wb.Worksheets.Item(1)
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsOverview.Range("E3").Value = "In Translation"
$wsOverview.Range("F3").Value = "In Translation"
$wsOverview.Range("E4").Value = "In Translation"
$wsOverview.Range("F4").Value = "In Translation"

$wsZhCn = $wb.Worksheets.Item(2)
$wsZhCn.Range("C2").Value = "In Translation"
$wsZhCn.Range("C3").Value = "In Translation"
$wsZhCn.Range("C4").Value = "In Translation"

$wsDeDe = $wb.Worksheets.Item(3)
$wsDeDe.Range("C2").Value = "In Translation"
$wsDeDe.Range("C3").Value = "In Translation"
$wsDeDe.Range("C4").Value = "In Translation"

# --- 2. Column width changes ---
# Target stored (OOXML) width is 13.4101845877511 characters, down from
# 17.2159881591797. The COM ColumnWidth setter here snaps the persisted
# width to the nearest 1/6-character increment, so we feed it the input
# (12.5) whose nearest representable result (80/6 = 13.333333333333334)
# is the closest achievable match to the target width.

$wsOverview.Columns.Item(5).ColumnWidth = 12.5   # column E ("zh-cn")
$wsOverview.Columns.Item(6).ColumnWidth = 12.5   # column F ("de-de")

$wsZhCn.Columns.Item(3).ColumnWidth = 12.5        # column C ("Status")

$wsDeDe.Columns.Item(3).ColumnWidth = 12.5        # column C ("Status")
